$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("cache size")
$ws1.Range("F40").Value = 999.999
$co = $ws1.ChartObjects().Item(1)
$chart = $co.Chart
$chart.SetSourceData($ws1.Range("A1:F54"))
